# Append a new data row (row 11) for the 2024-05-09 session, following the
# same layout as the existing good-trial rows (Date, then 6 grasp counts).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 20240509
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 5

# Mirror the saved selection/active cell state captured in the edit (A11).
$null = $ws.Range("A1").Select()
$null = $ws.Range("A11").Select()
